$wb = $excel.ActiveWorkbook

# ---- Sheet: runs ----
$ws = $wb.Worksheets.Item("runs")
$ws.Range("B2").Value = 11053.963500000002
$ws.Range("C2").Value = 8858.097
$ws.Range("D2").Value = 367
$ws.Range("E2").Value = 11173.613950000003
$ws.Range("F2").Value = 70
$ws.Range("H2").Value = 8858.097
$ws.Range("J2").Value = 1579801
$ws.Range("B3").Value = 10887.736
$ws.Range("C3").Value = 8794.002999999999
$ws.Range("D3").Value = 331
$ws.Range("E3").Value = 10997.693700000002
$ws.Range("F3").Value = 66.5
$ws.Range("H3").Value = 8794.002999999999
$ws.Range("J3").Value = 1580806
$ws.Range("B4").Value = 10940.431500000002
$ws.Range("C4").Value = 8822.076
$ws.Range("D4").Value = 326
$ws.Range("E4").Value = 11051.019250000001
$ws.Range("F4").Value = 68.5
$ws.Range("H4").Value = 8822.076
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 1579818
$ws.Range("B5").Value = 10940.431500000002
$ws.Range("C5").Value = 8822.076
$ws.Range("D5").Value = 331
$ws.Range("E5").Value = 11051.019250000001
$ws.Range("F5").Value = 68.5
$ws.Range("H5").Value = 8822.076
$ws.Range("J5").Value = 1579818

# ---- Sheet: per_resource_all_runs ----
$ws = $wb.Worksheets.Item("per_resource_all_runs")
$ws.Range("E2").Value = 17389
$ws.Range("G2").Value = 0.005056947256999999
$ws.Range("H2").Value = 0.006661402464599999
$ws.Range("E6").Value = 5323
$ws.Range("G6").Value = 0.0015479975989999997
$ws.Range("H6").Value = 0.0020391422922
$ws.Range("E7").Value = 3198
$ws.Range("G7").Value = 0.0009300199739999998
$ws.Range("H7").Value = 0.0012250943172000001
$ws.Range("E8").Value = 458194
$ws.Range("G8").Value = 0.13324877172199998
$ws.Range("H8").Value = 0.17552559899159997
$ws.Range("B17").Value = "https://www.youtube.com/s/player/ef5f17ca/www-widgetapi.vflset/www-widgetapi.js"
$ws.Range("E17").Value = 10433
$ws.Range("G17").Value = 0.0030340520289999996
$ws.Range("H17").Value = 0.0039966882462000005
$ws.Range("E19").Value = 13611
$ws.Range("G19").Value = 0.003958255742999999
$ws.Range("H19").Value = 0.0052141209354
$ws.Range("E26").Value = 23186
$ws.Range("G26").Value = 0.006742790217999999
$ws.Range("H26").Value = 0.008882125340400001
$ws.Range("E27").Value = 56254
$ws.Range("G27").Value = 0.016359394502
$ws.Range("H27").Value = 0.021549861075600002
$ws.Range("E28").Value = 251978
$ws.Range("G28").Value = 0.073278478114
$ws.Range("H28").Value = 0.0965280850092
$ws.Range("E29").Value = 19944
$ws.Range("G29").Value = 0.005799974471999999
$ws.Range("H29").Value = 0.0076401754416
$ws.Range("E30").Value = 22174
$ws.Range("F30").Value = 257200
$ws.Range("G30").Value = 0.006448487461999999
$ws.Range("H30").Value = 0.008494446963600001
$ws.Range("B32").Value = "https://6019370.global.siteimproveanalytics.io/image.aspx?url=https%3A%2F%2Fwww.ullensaker.kommune.no%2F&title=Ullensaker&res=412x823&accountid=6019370&rt=812&luid=764eb221-2690-4560-d0b9-66f668c11bf0&dnt=true&prev=a41cf251-1740-3470-c430-485e4ead84bc&rnd=22734"
$ws.Range("E37").Value = 7416
$ws.Range("G37").Value = 0.0021566692079999995
$ws.Range("H37").Value = 0.0028409316624000004
$ws.Range("E38").Value = 10400
$ws.Range("G38").Value = 0.0030244551999999997
$ws.Range("H38").Value = 0.003984046560000001
$ws.Range("E40").Value = 58
$ws.Range("G40").Value = 0.000016867153999999998
$ws.Range("H40").Value = 0.000022218721200000005
$ws.Range("E41").Value = 1265
$ws.Range("G41").Value = 0.000367878445
$ws.Range("H41").Value = 0.00048459797099999994
$ws.Range("E42").Value = 17388
$ws.Range("G42").Value = 0.005056656444
$ws.Range("H42").Value = 0.0066610193832
$ws.Range("E43").Value = 645
$ws.Range("G43").Value = 0.00018757438499999996
$ws.Range("H43").Value = 0.000247087503
$ws.Range("E45").Value = 110943
$ws.Range("G45").Value = 0.03226366665899999
$ws.Range("H45").Value = 0.0425001997602
$ws.Range("E46").Value = 5323
$ws.Range("G46").Value = 0.0015479975989999997
$ws.Range("H46").Value = 0.0020391422922
$ws.Range("E47").Value = 3577
$ws.Range("G47").Value = 0.0010402381009999998
$ws.Range("H47").Value = 0.0013702821678000003
$ws.Range("E49").Value = 54
$ws.Range("G49").Value = 0.000015703902
$ws.Range("H49").Value = 0.0000206863956
$ws.Range("B57").Value = "https://www.youtube.com/s/player/ef5f17ca/www-widgetapi.vflset/www-widgetapi.js"
$ws.Range("E58").Value = 3340
$ws.Range("G58").Value = 0.0009713154199999999
$ws.Range("H58").Value = 0.0012794918760000002
$ws.Range("E59").Value = 13625
$ws.Range("G59").Value = 0.003962327125
$ws.Range("H59").Value = 0.005219484075000001
$ws.Range("E66").Value = 23467
$ws.Range("G66").Value = 0.006824508670999999
$ws.Range("H66").Value = 0.0089897712138
$ws.Range("E67").Value = 56235
$ws.Range("G67").Value = 0.016353869054999998
$ws.Range("H67").Value = 0.021542582529
$ws.Range("E68").Value = 252106
$ws.Range("G68").Value = 0.073315702178
$ws.Range("H68").Value = 0.0965771194284
$ws.Range("E69").Value = 20166
$ws.Range("G69").Value = 0.005864534957999999
$ws.Range("H69").Value = 0.0077252195124
$ws.Range("E70").Value = 22174
$ws.Range("F70").Value = 257200
$ws.Range("G70").Value = 0.006448487461999999
$ws.Range("H70").Value = 0.008494446963600001
$ws.Range("B72").Value = "https://6019370.global.siteimproveanalytics.io/image.aspx?url=https%3A%2F%2Fwww.ullensaker.kommune.no%2F&title=Ullensaker&res=412x823&accountid=6019370&rt=749&luid=97134075-012b-dcf4-0d3f-4bd596da7556&dnt=true&prev=1b1f88f7-9917-5ebf-59c8-d1a3e355d290&rnd=53132"
$ws.Range("E74").Value = 2444
$ws.Range("G74").Value = 0.0007107469719999999
$ws.Range("H74").Value = 0.0009362509416
$ws.Range("E77").Value = 7417
$ws.Range("G77").Value = 0.0021569600209999998
$ws.Range("H77").Value = 0.0028413147438
$ws.Range("E81").Value = 1288
$ws.Range("G81").Value = 0.000374567144
$ws.Range("H81").Value = 0.0004934088432000001
$ws.Range("E82").Value = 17388
$ws.Range("G82").Value = 0.005056656444
$ws.Range("H82").Value = 0.0066610193832
$ws.Range("E83").Value = 645
$ws.Range("G83").Value = 0.00018757438499999996
$ws.Range("H83").Value = 0.000247087503
$ws.Range("E86").Value = 5323
$ws.Range("G86").Value = 0.0015479975989999997
$ws.Range("H86").Value = 0.0020391422922
$ws.Range("E87").Value = 3198
$ws.Range("G87").Value = 0.0009300199739999998
$ws.Range("H87").Value = 0.0012250943172000001
$ws.Range("B97").Value = "https://www.ullensaker.kommune.no/siteassets/30-bilder/logoer/ullensaker.png"
$ws.Range("C97").Value = "Image"
$ws.Range("D97").Value = "image/png"
$ws.Range("E97").Value = 3340
$ws.Range("F97").Value = 3251
$ws.Range("G97").Value = 0.0009713154199999999
$ws.Range("H97").Value = 0.0012794918760000002
$ws.Range("B98").Value = "https://www.youtube.com/s/player/ef5f17ca/www-widgetapi.vflset/www-widgetapi.js"
$ws.Range("C98").Value = "Script"
$ws.Range("D98").Value = "text/javascript"
$ws.Range("E98").Value = 10432
$ws.Range("F98").Value = 30432
$ws.Range("G98").Value = 0.0030337612159999994
$ws.Range("H98").Value = 0.0039963051648
$ws.Range("E99").Value = 13622
$ws.Range("G99").Value = 0.0039614546859999995
$ws.Range("H99").Value = 0.005218334830800001
$ws.Range("E107").Value = 56245
$ws.Range("G107").Value = 0.016356777184999997
$ws.Range("H107").Value = 0.021546413343000004
$ws.Range("E108").Value = 251985
$ws.Range("G108").Value = 0.07328051380499999
$ws.Range("H108").Value = 0.09653076657899999
$ws.Range("E109").Value = 19944
$ws.Range("G109").Value = 0.005799974471999999
$ws.Range("H109").Value = 0.0076401754416
$ws.Range("E110").Value = 22174
$ws.Range("F110").Value = 257200
$ws.Range("G110").Value = 0.006448487461999999
$ws.Range("H110").Value = 0.008494446963600001
$ws.Range("B112").Value = "https://6019370.global.siteimproveanalytics.io/image.aspx?url=https%3A%2F%2Fwww.ullensaker.kommune.no%2F&title=Ullensaker&res=412x823&accountid=6019370&rt=771&luid=c18ff15d-a7c2-4e20-a9f0-0f7c8c5b09d9&dnt=true&prev=ab54e41c-6f15-130e-ec6c-64fb88f29a2b&rnd=80341"
$ws.Range("E114").Value = 2444
$ws.Range("G114").Value = 0.0007107469719999999
$ws.Range("H114").Value = 0.0009362509416
$ws.Range("E120").Value = 58
$ws.Range("G120").Value = 0.000016867153999999998
$ws.Range("H120").Value = 0.000022218721200000005
$ws.Range("E121").Value = 1288
$ws.Range("G121").Value = 0.000374567144
$ws.Range("H121").Value = 0.0004934088432000001

# ---- Sheet: summary_by_type ----
$ws = $wb.Worksheets.Item("summary_by_type")
$ws.Range("C2").Value = 52165
$ws.Range("C3").Value = 357864
$ws.Range("C4").Value = 2616379
$ws.Range("C5").Value = 336
$ws.Range("C6").Value = 1150233
$ws.Range("C8").Value = 66522
$ws.Range("D8").Value = 771600
$ws.Range("C11").Value = 3841

# ---- Sheet: co2 ----
$ws = $wb.Worksheets.Item("co2")
$ws.Range("B2").Value = 1579801
$ws.Range("C2").Value = 0.459426668213
$ws.Range("D2").Value = 0.6051923788014001
$ws.Range("B3").Value = 1580806
$ws.Range("C3").Value = 0.45971893527799995
$ws.Range("D3").Value = 0.6055773756084002
$ws.Range("B4").Value = 1579818
$ws.Range("C4").Value = 0.4594316120339999
$ws.Range("D4").Value = 0.6051988911852001
$ws.Range("B5").Value = 1579818
$ws.Range("C5").Value = 0.4594316120339999
$ws.Range("D5").Value = 0.6051988911852001
